$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value for this update
$updates = [ordered]@{
    'D2' = '67.479.11'
    'E2' = '  -0.56%  '
    'D3' = '3.718.46'
    'E3' = '  -2.19%  '
    'E4' = '  -0.16%  '
    'D5' = '589.19'
    'E5' = '  -1.69%  '
    'D6' = '164.82'
    'E6' = '  -2.52%  '
    'D7' = '3.716.53'
    'E7' = '  -2.22%  '
    'E8' = '  -0.02%  '
    'D9' = '0.516'
    'E9' = '  -1.81%  '
    'E10' = '  -4.16%  '
    'D11' = '6.40'
    'E11' = '  -1.16%  '
    'D12' = '0.448'
    'E12' = '  -2.27%  '
    'D13' = '0.0000260'
    'E13' = '  -6.52%  '
    'D14' = '35.77'
    'E14' = '  -2.61%  '
    'D15' = '4.342.51'
    'E15' = '  -2.15%  '
    'D16' = '3.718.98'
    'E16' = '  -1.91%  '
    'D17' = '67.398.26'
    'E17' = '  -0.88%  '
    'D18' = '18.25'
    'E18' = '  +0.88%  '
    'D19' = '7.01'
    'E19' = '  -5.11%  '
    'E20' = '  -0.21%  '
    'D21' = '10.58'
    'E21' = '  -2.06%  '
    'D22' = '463.52'
    'E22' = '  -0.92%  '
    'D23' = '0.697'
    'E23' = '  -3.78%  '
    'D24' = '82.26'
    'E24' = '  -1.08%  '
    'E25' = '  -11.36%  '
    'B26' = 'Fetch.AI'
    'C26' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'D26' = '2.16'
    'E26' = '  -3.84%  '
    'B27' = 'InternetComputer(DFINITY)'
    'C27' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D27' = '11.91'
    'E27' = '  -1.49%  '
    'D28' = '10.16'
    'E28' = '  -0.93%  '
    'E29' = '  +0.01%  '
    'D30' = '3.862.33'
    'E30' = '  -2.18%  '
    'D31' = '2.75'
    'E31' = '  -6.15%  '
    'D32' = '7.29'
    'E32' = '  -5.29%  '
    'D33' = '2.19'
    'E33' = '  -3.81%  '
    'D34' = '29.56'
    'E34' = '  -3.76%  '
    'D35' = '8.96'
    'E35' = '  -3.87%  '
    'D36' = '3.671.19'
    'E36' = '  -2.53%  '
    'D37' = '0.101'
    'E37' = '  -5.22%  '
    'D38' = '3.40'
    'E38' = '  -10.92%  '
    'D39' = '0.987'
    'E39' = '  -2.48%  '
    'D40' = '0.135'
    'E40' = '  -3.32%  '
    'D41' = '5.71'
    'E41' = '  -3.72%  '
    'D42' = '0.998'
    'E42' = '  -0.22%  '
    'E43' = '  -0.01%  '
    'D44' = '0.303'
    'E44' = '  -3.72%  '
    'D45' = '8.48'
    'E45' = '  -3.66%  '
    'D46' = '1.91'
    'E46' = '  -3.06%  '
    'D47' = '45.19'
    'E47' = '  -2.78%  '
    'D48' = '392.15'
    'E48' = '  -4.04%  '
    'D49' = '143.30'
    'E49' = '  +1.39%  '
    'D50' = '0.0345'
    'E50' = '  -3.65%  '
    'D51' = '24.99'
    'E51' = '  -2.34%  '
}

# Cells whose new value looks numeric but must stay stored as literal text
# (matches source data which keeps these as inline strings, not numbers)
$textForceCells = @('D5', 'D6', 'D9', 'D11', 'D12', 'D13', 'D14', 'D18', 'D19', 'D21', 'D22', 'D23', 'D24', 'D26', 'D27', 'D28', 'D31', 'D32', 'D33', 'D34', 'D35', 'D37', 'D38', 'D39', 'D40', 'D41', 'D42', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')

foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}

Write-Host "Applied $($updates.Count) cell updates"
